$wb = $excel.ActiveWorkbook

# Add violent crime data for 2023-12-25, distributed across
# citywide totals (by crime category), the by-neighborhood summary,
# and each affected neighborhood own sheet (by crime category).

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 7532  # Aggravated Assault: 7513 -> 7532
$ws.Range('J3').Value = 7936  # Aggravated Battery: 7903 -> 7936
$ws.Range('J4').Value = 1727  # Criminal Sexual Assault: 1721 -> 1727
$ws.Range('J6').Value = 10834  # Robbery: 10802 -> 10834
$ws.Range('J7').Value = 28648  # Total: 28558 -> 28648

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('J2').Value = 22  # Aggravated Assault: 21 -> 22
$ws.Range('J7').Value = 55  # Total: 54 -> 55

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J3').Value = 57  # Aggravated Battery: 56 -> 57
$ws.Range('J7').Value = 428  # Total: 427 -> 428

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 479  # Aggravated Assault: 477 -> 479
$ws.Range('J3').Value = 519  # Aggravated Battery: 517 -> 519
$ws.Range('J6').Value = 670  # Robbery: 669 -> 670
$ws.Range('J7').Value = 1806  # Total: 1801 -> 1806

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 293  # Aggravated Assault: 292 -> 293
$ws.Range('J3').Value = 431  # Aggravated Battery: 430 -> 431
$ws.Range('J6').Value = 462  # Robbery: 461 -> 462
$ws.Range('J7').Value = 1296  # Total: 1293 -> 1296

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 295  # Aggravated Battery: 292 -> 295
$ws.Range('J6').Value = 256  # Robbery: 254 -> 256
$ws.Range('J7').Value = 878  # Total: 873 -> 878

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 209  # Aggravated Assault: 208 -> 209
$ws.Range('J6').Value = 268  # Robbery: 267 -> 268
$ws.Range('J7').Value = 718  # Total: 716 -> 718

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 172  # Aggravated Battery: 168 -> 172
$ws.Range('J7').Value = 436  # Total: 432 -> 436

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J4').Value = 131  # Archer Heights: 130 -> 131
$ws.Range('J7').Value = 819  # Auburn Gresham: 816 -> 819
$ws.Range('J8').Value = 1806  # Austin: 1801 -> 1806
$ws.Range('J11').Value = 522  # Belmont Cragin: 516 -> 522
$ws.Range('J15').Value = 355  # Brighton Park: 353 -> 355
$ws.Range('J19').Value = 835  # Chatham: 827 -> 835
$ws.Range('J20').Value = 621  # Chicago Lawn: 617 -> 621
$ws.Range('J21').Value = 82  # Chinatown: 80 -> 82
$ws.Range('J24').Value = 102  # Dunning: 101 -> 102
$ws.Range('J27').Value = 174  # Edgewater: 173 -> 174
$ws.Range('J29').Value = 1529  # Englewood: 1525 -> 1529
$ws.Range('J31').Value = 301  # Gage Park: 299 -> 301
$ws.Range('J33').Value = 1296  # Garfield Park: 1293 -> 1296
$ws.Range('J34').Value = 131  # Garfield Ridge: 130 -> 131
$ws.Range('J36').Value = 389  # Grand Boulevard: 387 -> 389
$ws.Range('J37').Value = 878  # Grand Crossing: 873 -> 878
$ws.Range('J42').Value = 1211  # Humboldt Park: 1210 -> 1211
$ws.Range('J48').Value = 318  # Lake View: 317 -> 318
$ws.Range('J50').Value = 178  # Lincoln Square: 176 -> 178
$ws.Range('J51').Value = 360  # Little Italy, UIC: 359 -> 360
$ws.Range('J52').Value = 727  # Little Village: 725 -> 727
$ws.Range('J53').Value = 428  # Logan Square: 427 -> 428
$ws.Range('J55').Value = 451  # Lower West Side: 450 -> 451
$ws.Range('J57').Value = 137  # Mckinley Park: 136 -> 137
$ws.Range('J63').Value = 86  # NO NEIGHBORHOOD DATA: 83 -> 86
$ws.Range('J65').Value = 718  # New City: 716 -> 718
$ws.Range('J67').Value = 1043  # North Lawndale: 1040 -> 1043
$ws.Range('J69').Value = 55  # Norwood Park: 54 -> 55
$ws.Range('J72').Value = 109  # Old Town: 108 -> 109
$ws.Range('J73').Value = 282  # Portage Park: 280 -> 282
$ws.Range('J77').Value = 201  # Riverdale: 200 -> 201
$ws.Range('J78').Value = 335  # Rogers Park: 333 -> 335
$ws.Range('J79').Value = 786  # Roseland: 784 -> 786
$ws.Range('J85').Value = 1177  # South Shore: 1175 -> 1177
$ws.Range('J86').Value = 173  # Streeterville: 172 -> 173
$ws.Range('J89').Value = 358  # Uptown: 356 -> 358
$ws.Range('J92').Value = 94  # West Elsdon: 92 -> 94
$ws.Range('J96').Value = 324  # West Ridge: 323 -> 324
$ws.Range('J98').Value = 209  # Wicker Park: 208 -> 209
$ws.Range('J99').Value = 436  # Woodlawn: 432 -> 436
$ws.Range('J101').Value = 28648  # Total: 28558 -> 28648

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J6').Value = 109  # Robbery: 107 -> 109
$ws.Range('J7').Value = 301  # Total: 299 -> 301

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 395  # Aggravated Battery: 392 -> 395
$ws.Range('J7').Value = 1043  # Total: 1040 -> 1043

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J3').Value = 114  # Aggravated Battery: 113 -> 114
$ws.Range('J6').Value = 257  # Robbery: 258 -> 257

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 463  # Aggravated Assault: 461 -> 463
$ws.Range('J3').Value = 537  # Aggravated Battery: 536 -> 537
$ws.Range('J6').Value = 390  # Robbery: 389 -> 390
$ws.Range('J7').Value = 1529  # Total: 1525 -> 1529

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J6').Value = 154  # Robbery: 153 -> 154
$ws.Range('J7').Value = 318  # Total: 317 -> 318

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 206  # Aggravated Assault: 204 -> 206
$ws.Range('J4').Value = 41  # Criminal Sexual Assault: 40 -> 41
$ws.Range('J6').Value = 326  # Robbery: 321 -> 326
$ws.Range('J7').Value = 835  # Total: 827 -> 835

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J6').Value = 641  # Robbery: 640 -> 641
$ws.Range('J7').Value = 1211  # Total: 1210 -> 1211

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J3').Value = 103  # Aggravated Battery: 101 -> 103
$ws.Range('J7').Value = 335  # Total: 333 -> 335

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J6').Value = 254  # Robbery: 253 -> 254
$ws.Range('J7').Value = 451  # Total: 450 -> 451

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J6').Value = 28  # Robbery: 27 -> 28
$ws.Range('J7').Value = 102  # Total: 101 -> 102

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J3').Value = 83  # Aggravated Battery: 82 -> 83
$ws.Range('J7').Value = 324  # Total: 323 -> 324

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('J2').Value = 14  # Aggravated Assault: 13 -> 14
$ws.Range('J6').Value = 55  # Robbery: 54 -> 55
$ws.Range('J7').Value = 82  # Total: 80 -> 82

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 227  # Aggravated Assault: 226 -> 227
$ws.Range('J3').Value = 262  # Aggravated Battery: 261 -> 262
$ws.Range('J7').Value = 786  # Total: 784 -> 786

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J3').Value = 204  # Aggravated Battery: 201 -> 204
$ws.Range('J6').Value = 181  # Robbery: 180 -> 181
$ws.Range('J7').Value = 621  # Total: 617 -> 621

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J6').Value = 116  # Robbery: 114 -> 116
$ws.Range('J7').Value = 389  # Total: 387 -> 389

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 257  # Aggravated Assault: 256 -> 257
$ws.Range('J3').Value = 247  # Aggravated Battery: 245 -> 247
$ws.Range('J7').Value = 819  # Total: 816 -> 819

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J3').Value = 35  # Aggravated Battery: 34 -> 35
$ws.Range('J7').Value = 131  # Total: 130 -> 131

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 74  # Aggravated Battery: 73 -> 74
$ws.Range('J6').Value = 163  # Robbery: 162 -> 163
$ws.Range('J7').Value = 355  # Total: 353 -> 355

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J3').Value = 30  # Aggravated Battery: 29 -> 30
$ws.Range('J7').Value = 209  # Total: 208 -> 209

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J2').Value = 46  # Aggravated Assault: 45 -> 46
$ws.Range('J6').Value = 60  # Robbery: 59 -> 60
$ws.Range('J7').Value = 178  # Total: 176 -> 178

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 142  # Aggravated Assault: 141 -> 142
$ws.Range('J6').Value = 252  # Robbery: 247 -> 252
$ws.Range('J7').Value = 522  # Total: 516 -> 522

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J6').Value = 104  # Robbery: 102 -> 104
$ws.Range('J7').Value = 282  # Total: 280 -> 282

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('J2').Value = 27  # Aggravated Assault: 26 -> 27
$ws.Range('J3').Value = 28  # Aggravated Battery: 27 -> 28
$ws.Range('J7').Value = 94  # Total: 92 -> 94

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J3').Value = 102  # Aggravated Battery: 101 -> 102
$ws.Range('J4').Value = 36  # Criminal Sexual Assault: 35 -> 36
$ws.Range('J7').Value = 358  # Total: 356 -> 358

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J2').Value = 46  # Aggravated Assault: 45 -> 46
$ws.Range('J7').Value = 174  # Total: 173 -> 174

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J4').Value = 94  # Criminal Sexual Assault: 93 -> 94
$ws.Range('J7').Value = 173  # Total: 172 -> 173

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J2').Value = 76  # Aggravated Assault: 75 -> 76
$ws.Range('J7').Value = 360  # Total: 359 -> 360

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J3').Value = 32  # Aggravated Battery: 31 -> 32
$ws.Range('J7').Value = 137  # Total: 136 -> 137

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 427  # Aggravated Battery: 426 -> 427
$ws.Range('J6').Value = 336  # Robbery: 335 -> 336
$ws.Range('J7').Value = 1177  # Total: 1175 -> 1177

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('J6').Value = 41  # Robbery: 40 -> 41
$ws.Range('J7').Value = 109  # Total: 108 -> 109

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J6').Value = 38  # Robbery: 37 -> 38
$ws.Range('J7').Value = 201  # Total: 200 -> 201

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 172  # Aggravated Assault: 171 -> 172
$ws.Range('J3').Value = 204  # Aggravated Battery: 203 -> 204
$ws.Range('J7').Value = 727  # Total: 725 -> 727

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J2').Value = 41  # Aggravated Assault: 40 -> 41
$ws.Range('J7').Value = 131  # Total: 130 -> 131
